# glottocreate: add/modify metatables for description, references and contributors. Fixes #49
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "metadata" -> "description" : rename + update header row
# ---------------------------------------------------------------------------
$description = $wb.Worksheets.Item("metadata")
$description.Name = "description"

$description.Range("D1").Value = "remarks"
$description.Range("E1").Value = "lev_Y"
$description.Range("F1").Value = "lev_N"
$description.Range("G1").Value = "lev_NA"
$description.Range("H1").Value = "lev_A"
$description.Range("I1").Value = "lev_B"
$description.Range("J1").Value = "lev_C"
$description.Range("K1").Value = "lev_D"

# ---------------------------------------------------------------------------
# 2) "references" : insert "reference" + "page" columns after glottocode
# ---------------------------------------------------------------------------
$references = $wb.Worksheets.Item("references")
$references.Range("B1:C1").Insert(-4161)
$references.Range("B1").Value = "reference"
$references.Range("C1").Value = "page"

# ---------------------------------------------------------------------------
# 3) "remarks" : insert "remark" column after glottocode
# ---------------------------------------------------------------------------
$remarks = $wb.Worksheets.Item("remarks")
$remarks.Range("B1").Insert(-4161)
$remarks.Range("B1").Value = "remark"

# ---------------------------------------------------------------------------
# 4) Recreate "readme" and "lookup" after a new "contributors" sheet so that
#    sheet order becomes: ... remarks, contributors, readme, lookup
# ---------------------------------------------------------------------------
$readmeOld = $wb.Worksheets.Item("readme")
$readmeInfo = @()
for ($r = 1; $r -le 6; $r++) {
    $a = $readmeOld.Range("A$r").Value2
    $b = $readmeOld.Range("B$r").Value2
    $readmeInfo += , @($a, $b)
}

$lookupOld = $wb.Worksheets.Item("lookup")
$lookupInfo = @()
for ($r = 1; $r -le 8; $r++) {
    $a = $lookupOld.Range("A$r").Value2
    $b = $lookupOld.Range("B$r").Value2
    $lookupInfo += , @($a, $b)
}

$wb.Worksheets.Item("readme").Delete()
$wb.Worksheets.Item("lookup").Delete()

$remarksSheet = $wb.Worksheets.Item("remarks")
$contributors = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $remarksSheet)
$contributors.Name = "contributors"

$contributors.Range("A1").Value = "glottocode"
$contributors.Range("B1").Value = "contributor"
$contributors.Range("C1").Value = "var001_contributor"
$contributors.Range("D1").Value = "var002_contributor"
$contributors.Range("E1").Value = "var003_contributor"
$contributors.Range("A2").Value = "yucu1253"
$contributors.Range("A3").Value = "tani1257"

$readme = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $contributors)
$readme.Name = "readme"
for ($r = 1; $r -le 6; $r++) {
    $pair = $readmeInfo[$r - 1]
    if ($pair[0] -ne $null) { $readme.Range("A$r").Value = $pair[0] }
    if ($pair[1] -ne $null) { $readme.Range("B$r").Value = $pair[1] }
}

$lookup = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $readme)
$lookup.Name = "lookup"
for ($r = 1; $r -le 8; $r++) {
    $pair = $lookupInfo[$r - 1]
    if ($pair[0] -ne $null) { $lookup.Range("A$r").Value = $pair[0] }
    if ($pair[1] -ne $null) { $lookup.Range("B$r").Value = $pair[1] }
}
